$d = $word.ActiveDocument

function Get-ParagraphByText($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Text.StartsWith($prefix)) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. Merge the "Open transform.py file" paragraph into the following
#    paragraph ("Go to line 1006 ...") by deleting the paragraph mark that
#    ends it - this is exactly what happens in Word when you press Delete at
#    the end of a paragraph: the two paragraphs become one, keeping the
#    paragraph mark (and therefore the pPr) of the second paragraph.
# ---------------------------------------------------------------------------
$openPara = Get-ParagraphByText "Open transform.py file"
$markRange = $d.Range($openPara.Range.End - 1, $openPara.Range.End)
$markRange.Delete()

# ---------------------------------------------------------------------------
# 2. The merged paragraph now reads:
#      "Open transform.py fileGo to line 1006 (this should be describing
#       keyword arguments for splitPixel.fullspit2D())change pos0Range
#       (Wingdings arrow) pos0_range"
#    Replace its whole text with the corrected wording, keeping the
#    paragraph's own formatting (pPr) untouched, and keeping the trailing
#    "change pos0Range -> pos0_range" line (with its manual line break and
#    Wingdings arrow symbol) exactly as it already is.
# ---------------------------------------------------------------------------
$mergedPara = Get-ParagraphByText "Open transform.py file"
$bodyRange = $d.Range($mergedPara.Range.Start, $mergedPara.Range.End - 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>Go to line 10</w:t></w:r>' + `
    '<w:r><w:t>32</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (this should be describing keyword arguments for split</w:t></w:r>' + `
    '<w:r><w:t>BBox</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
    '<w:r><w:t>histoBBox2</w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>d</w:t></w:r>' + `
    '<w:r><w:t>(</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t>))</w:t></w:r>' + `
    '<w:r><w:br/><w:t xml:space="preserve">change pos0Range </w:t></w:r>' + `
    '<w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' + `
    '<w:r><w:t xml:space="preserve"> pos0_range</w:t></w:r>' + `
    '</w:p>'
$bodyRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3. Insert a new, empty paragraph (matching the surrounding indentation)
#    right after the "change pos1Range -> pos1_range" paragraph.
# ---------------------------------------------------------------------------
$pos1Para = Get-ParagraphByText "change pos1Range"
$insertPoint = $d.Range($pos1Para.Range.End, $pos1Para.Range.End)
$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/></w:pPr></w:p>'
$null = $insertPoint.InsertXML($emptyParaXml)
